$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996

$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984

$ws.Range("H86").Value = 1430.4
$ws.Range("I86").Value = 1400
$ws.Range("J86").Value = 1450.6666
$ws.Range("K86").Value = 1400
$ws.Range("L86").Value = 1450.6666
$ws.Range("M86").Value = -277
$ws.Range("N86").Value = -3696.6666

$ws.Range("H89").Value = 1430.4
$ws.Range("I89").Value = 1400
$ws.Range("J89").Value = 1450.6666
$ws.Range("K89").Value = 7000
$ws.Range("L89").Value = 7253.333000000001
$ws.Range("M89").Value = -1384
$ws.Range("N89").Value = -18485.333

$ws.Range("H113").Value = 12573.048
$ws.Range("I113").Value = 23598.223
$ws.Range("J113").Value = 4304.1665
$ws.Range("K113").Value = 23598.223
$ws.Range("L113").Value = 4304.1665
$ws.Range("M113").Value = -20344.223
$ws.Range("N113").Value = -10812.1665

$ws.Range("H132").Value = 34960584
$ws.Range("I132").Value = 37720228
$ws.Range("J132").Value = 5100
$ws.Range("K132").Value = 113160684
$ws.Range("L132").Value = 15300
$ws.Range("M132").Value = -113158154
$ws.Range("N132").Value = -20360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3542.309
$ws.Range("I32").Value = 2798.0667
$ws.Range("J32").Value = 6891.4
$ws.Range("K32").Value = 2798.0667
$ws.Range("L32").Value = 6891.4
$ws.Range("M32").Value = -2511.0667
$ws.Range("N32").Value = -7465.4

$ws.Range("H74").Value = 5927.7617
$ws.Range("I74").Value = 573.46155
$ws.Range("K74").Value = 573.46155
$ws.Range("M74").Value = 300.53845

$ws.Range("H77").Value = 5927.7617
$ws.Range("I77").Value = 573.46155
$ws.Range("K77").Value = 2867.30775
$ws.Range("M77").Value = 1500.69225

$ws.Range("H122").Value = 63158680
$ws.Range("I122").Value = 70589000
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 211767000
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -211764550
$ws.Range("N122").Value = -7900

$ws.Range("H135").Value = 129809.664
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 129809.664
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 129809.664
$ws.Range("M135").Value = $null
$ws.Range("N135").Value = -139949.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8560632
$ws.Range("I134").Value = 12364814
$ws.Range("J134").Value = 1221
$ws.Range("K134").Value = 37094442
$ws.Range("L134").Value = 3663
$ws.Range("M134").Value = -37091907
$ws.Range("N134").Value = -8733

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2718806.5
$ws.Range("I58").Value = 3509332
$ws.Range("J58").Value = 17843.666
$ws.Range("K58").Value = 3509332
$ws.Range("L58").Value = 17843.666
$ws.Range("M58").Value = -3509129
$ws.Range("N58").Value = -18249.666

$ws.Range("H86").Value = 3297.1
$ws.Range("I86").Value = 3162.8667
$ws.Range("J86").Value = 3699.8
$ws.Range("K86").Value = 3162.8667
$ws.Range("L86").Value = 3699.8
$ws.Range("M86").Value = -2039.8667
$ws.Range("N86").Value = -5945.8

$ws.Range("H89").Value = 3297.1
$ws.Range("I89").Value = 3162.8667
$ws.Range("J89").Value = 3699.8
$ws.Range("K89").Value = 15814.3335
$ws.Range("L89").Value = 18499
$ws.Range("M89").Value = -10198.3335
$ws.Range("N89").Value = -29731

$ws.Range("H132").Value = 9263487
$ws.Range("I132").Value = 13889868
$ws.Range("J132").Value = 10726.167
$ws.Range("K132").Value = 41669604
$ws.Range("L132").Value = 32178.501
$ws.Range("M132").Value = -41667074
$ws.Range("N132").Value = -37238.501

$ws.Range("H134").Value = 7268385.5
$ws.Range("I134").Value = 8334066.5
$ws.Range("J134").Value = 4809122
$ws.Range("K134").Value = 25002199.5
$ws.Range("L134").Value = 14427366
$ws.Range("M134").Value = -24999664.5
$ws.Range("N134").Value = -14432436

$ws.Range("H136").Value = 2718806.5
$ws.Range("I136").Value = 3509332
$ws.Range("J136").Value = 17843.666
$ws.Range("K136").Value = 10527996
$ws.Range("L136").Value = 53530.99800000001
$ws.Range("M136").Value = -10525446
$ws.Range("N136").Value = -58630.99800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 723.2222
$ws.Range("I5").Value = 707.86664
$ws.Range("J5").Value = 800
$ws.Range("K5").Value = 2123.59992
$ws.Range("L5").Value = 2400
$ws.Range("M5").Value = -2011.59992
$ws.Range("N5").Value = -2624

$ws.Range("H97").Value = 1200
$ws.Range("I97").Value = 750
$ws.Range("J97").Value = 1380
$ws.Range("K97").Value = 2250
$ws.Range("L97").Value = 4140
$ws.Range("M97").Value = -1754
$ws.Range("N97").Value = -5132

$ws.Range("H98").Value = 5432.6665
$ws.Range("I98").Value = 4500
$ws.Range("J98").Value = 5699.143
$ws.Range("K98").Value = 13500
$ws.Range("L98").Value = 17097.429
$ws.Range("M98").Value = -12002
$ws.Range("N98").Value = -20093.429

$ws.Range("H107").Value = 7784.5625
$ws.Range("I107").Value = 275.375
$ws.Range("J107").Value = 15293.75
$ws.Range("K107").Value = 826.125
$ws.Range("L107").Value = 45881.25
$ws.Range("M107").Value = 1093.875
$ws.Range("N107").Value = -49721.25

$ws.Range("H132").Value = 4000.5715
$ws.Range("I132").Value = 665.2381
$ws.Range("J132").Value = 9003.571
$ws.Range("K132").Value = 5987.142900000001
$ws.Range("L132").Value = 81032.139
$ws.Range("M132").Value = -3457.142900000001
$ws.Range("N132").Value = -86092.139

$ws.Range("H135").Value = 723.2222
$ws.Range("I135").Value = 707.86664
$ws.Range("J135").Value = 800
$ws.Range("K135").Value = 6370.79976
$ws.Range("L135").Value = 7200
$ws.Range("M135").Value = -3835.79976
$ws.Range("N135").Value = -12270

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2879.4736
$ws.Range("I80").Value = 2478.5715
$ws.Range("J80").Value = 3113.3333
$ws.Range("K80").Value = 2478.5715
$ws.Range("L80").Value = 3113.3333
$ws.Range("M80").Value = -1480.5715
$ws.Range("N80").Value = -5109.3333

$ws.Range("H83").Value = 2879.4736
$ws.Range("I83").Value = 2478.5715
$ws.Range("J83").Value = 3113.3333
$ws.Range("K83").Value = 12392.8575
$ws.Range("L83").Value = 15566.6665
$ws.Range("M83").Value = -7400.8575
$ws.Range("N83").Value = -25550.6665

$ws.Range("H134").Value = 30114.285
$ws.Range("J134").Value = 30114.285
$ws.Range("L134").Value = 90342.855
$ws.Range("N134").Value = -95412.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 4947.5
$ws.Range("J43").Value = 4947.5
$ws.Range("L43").Value = 4947.5
$ws.Range("N43").Value = -5333.5

$ws.Range("H61").Value = 2841
$ws.Range("I61").Value = 2123.5557
$ws.Range("J61").Value = 3428
$ws.Range("K61").Value = 2123.5557
$ws.Range("L61").Value = 3428
$ws.Range("M61").Value = -1921.5557
$ws.Range("N61").Value = -3832

$ws.Range("H82").Value = 4315
$ws.Range("I82").Value = 2250
$ws.Range("J82").Value = 4831.25
$ws.Range("K82").Value = 2250
$ws.Range("L82").Value = 4831.25
$ws.Range("M82").Value = -1889
$ws.Range("N82").Value = -5553.25

$ws.Range("H85").Value = 4315
$ws.Range("I85").Value = 2250
$ws.Range("J85").Value = 4831.25
$ws.Range("K85").Value = 2250
$ws.Range("L85").Value = 4831.25
$ws.Range("M85").Value = -1002
$ws.Range("N85").Value = -7327.25

$ws.Range("H113").Value = 2841
$ws.Range("I113").Value = 2123.5557
$ws.Range("J113").Value = 3428
$ws.Range("K113").Value = 2123.5557
$ws.Range("L113").Value = 3428
$ws.Range("M113").Value = 46.44430000000011
$ws.Range("N113").Value = -7768

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 359
$ws.Range("I100").Value = 182
$ws.Range("J100").Value = 801.5
$ws.Range("K100").Value = 364
$ws.Range("L100").Value = 1603
$ws.Range("M100").Value = 177
$ws.Range("N100").Value = -2685

$ws.Range("H126").Value = 250000670
$ws.Range("I126").Value = 156250750
$ws.Range("K126").Value = 468752250
$ws.Range("M126").Value = -468749780
